$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quadratic-svm-score")

# The commit "updated single child outputs" removes the stale
# "RUG357.fasta" prediction row (worksheet row 2) from the results
# table; Excel shifts every following row up by one and the
# dimension/shared-strings tables are recomputed automatically.
$ws.Rows.Item(2).Delete()
